$wb = $excel.ActiveWorkbook

# --- Sheet: Metadata ---
$meta = $wb.Worksheets.Item("Metadata")

# Date property (row 8, column B)
$meta.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# FHIR Version property (row 15, column B)
$meta.Range("B15").Value = "4.0.1"

# --- Sheet: Elements ---
$elements = $wb.Worksheets.Item("Elements")

# Extension.id (row 3) Type(s) column K: "id" -> "string"
$elements.Range("K3").Value = "string`n"

# Extension (row 2) Constraint(s) column AJ: drop the "unless an empty
# Parameters resource ... or $this is Parameters" clause from ele-1
$elements.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Extension.extension (row 4) Constraint(s) column AJ: same trimmed text
$elements.Range("AJ4").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Extension.value[x] (row 6) Definition column M: R4B -> R4 link
$elements.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."

# Extension.value[x]:valueQuantity (row 7) Definition column M: R4B -> R4 link
$elements.Range("M7").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
